$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.30259999999999
$ws.Range("C6").Value = -11.9147
$ws.Range("C7").Value = -12.7844
$ws.Range("C16").Value = -13.9884
$ws.Range("C20").Value = -12.28509999999999
$ws.Range("C28").Value = -13.09130000000001
$ws.Range("C29").Value = -11.79330000000001
$ws.Range("C32").Value = -12.14739999999999
$ws.Range("C40").Value = -12.35260000000001
$ws.Range("C46").Value = -14.27449999999998
$ws.Range("C51").Value = -11.8144
$ws.Range("C52").Value = -11.08450000000001
$ws.Range("C57").Value = -14.2274
$ws.Range("C59").Value = -12.4337
$ws.Range("C62").Value = -14.19169999999999
$ws.Range("C66").Value = -11.3232
$ws.Range("C73").Value = -11.03970000000001
$ws.Range("C74").Value = -12.0106
$ws.Range("C92").Value = -10.6269
$ws.Range("C100").Value = -11.1162
